$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 158
$ws.Range("I9").Value = 64.8
$ws.Range("J9").Value = 313.33334
$ws.Range("K9").Value = 64.8
$ws.Range("L9").Value = 313.33334
$ws.Range("M9").Value = 104.2
$ws.Range("N9").Value = -651.33334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 849.3484999999999
$ws.Range("J129").Value = 859.0323
$ws.Range("L129").Value = 2577.0969
$ws.Range("N129").Value = -12577.0969

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3201.0833
$ws.Range("I132").Value = 3401.6365
$ws.Range("J132").Value = 995
$ws.Range("K132").Value = 10204.9095
$ws.Range("L132").Value = 2985
$ws.Range("M132").Value = -7674.9095
$ws.Range("N132").Value = -8045

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 43626.668
$ws.Range("I137").Value = 1952.7222
$ws.Range("K137").Value = 5858.1666
$ws.Range("M137").Value = -3308.1666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2506.3333
$ws.Range("I138").Value = 1463.1364
$ws.Range("J138").Value = 3653.85
$ws.Range("K138").Value = 4389.4092
$ws.Range("L138").Value = 10961.55
$ws.Range("M138").Value = 750.5907999999999
$ws.Range("N138").Value = -21241.55

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3747.3572
$ws.Range("I45").Value = 3891.2
$ws.Range("J45").Value = 3667.4443
$ws.Range("K45").Value = 3891.2
$ws.Range("L45").Value = 3667.4443
$ws.Range("M45").Value = -3514.2
$ws.Range("N45").Value = -4421.4443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2545.9375
$ws.Range("I74").Value = 2101.6155
$ws.Range("J74").Value = 4471.3335
$ws.Range("K74").Value = 2101.6155
$ws.Range("L74").Value = 4471.3335
$ws.Range("M74").Value = -1227.6155
$ws.Range("N74").Value = -6219.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2545.9375
$ws.Range("I77").Value = 2101.6155
$ws.Range("J77").Value = 4471.3335
$ws.Range("K77").Value = 10508.0775
$ws.Range("L77").Value = 22356.6675
$ws.Range("M77").Value = -6140.077499999999
$ws.Range("N77").Value = -31092.6675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1982.4445
$ws.Range("I122").Value = 2105.25
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 6315.75
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -3865.75
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2106.8125
$ws.Range("I99").Value = 1818.25
$ws.Range("J99").Value = 2972.5
$ws.Range("K99").Value = 1818.25
$ws.Range("L99").Value = 2972.5
$ws.Range("M99").Value = -320.25
$ws.Range("N99").Value = -5968.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14853.5
$ws.Range("I31").Value = 19214.857
$ws.Range("K31").Value = 19214.857
$ws.Range("M31").Value = -18919.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 14853.5
$ws.Range("I34").Value = 19214.857
$ws.Range("K34").Value = 19214.857
$ws.Range("M34").Value = -19012.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 68663.336
$ws.Range("J68").Value = 68663.336
$ws.Range("L68").Value = 68663.336
$ws.Range("N68").Value = -70161.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 68663.336
$ws.Range("J71").Value = 68663.336
$ws.Range("L71").Value = 205990.008
$ws.Range("N71").Value = -213478.008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 16487.092
$ws.Range("I86").Value = 4083.3333
$ws.Range("J86").Value = 31371.6
$ws.Range("K86").Value = 4083.3333
$ws.Range("L86").Value = 31371.6
$ws.Range("M86").Value = -2960.3333
$ws.Range("N86").Value = -33617.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 16487.092
$ws.Range("I89").Value = 4083.3333
$ws.Range("J89").Value = 31371.6
$ws.Range("K89").Value = 20416.6665
$ws.Range("L89").Value = 156858
$ws.Range("M89").Value = -14800.6665
$ws.Range("N89").Value = -168090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 38000
$ws.Range("J133").Value = 38000
$ws.Range("L133").Value = 38000
$ws.Range("N133").Value = -43060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 930.7368
$ws.Range("I134").Value = 804.94116
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2414.82348
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = 120.17652
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 50422.5
$ws.Range("J135").Value = 50422.5
$ws.Range("L135").Value = 50422.5
$ws.Range("N135").Value = -60562.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 137.33333
$ws.Range("J15").Value = 201
$ws.Range("L15").Value = 603
$ws.Range("N15").Value = -883

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 221.25
$ws.Range("I33").Value = 195
$ws.Range("K33").Value = 1170
$ws.Range("M33").Value = -887

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1500
$ws.Range("J80").Value = 1500
$ws.Range("L80").Value = 4500
$ws.Range("N80").Value = -6372

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 1500
$ws.Range("J83").Value = 1500
$ws.Range("L83").Value = 13500
$ws.Range("N83").Value = -22860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 491.33334
$ws.Range("J122").Value = 491.33334
$ws.Range("L122").Value = 4422.00006
$ws.Range("N122").Value = -9322.00006

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 101811.96
$ws.Range("J131").Value = 104962.02
$ws.Range("L131").Value = 314886.06
$ws.Range("N131").Value = -324966.06

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4773.25
$ws.Range("I70").Value = 4745
$ws.Range("J70").Value = 4787.375
$ws.Range("K70").Value = 4745
$ws.Range("L70").Value = 4787.375
$ws.Range("M70").Value = -4475
$ws.Range("N70").Value = -5327.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4773.25
$ws.Range("I73").Value = 4745
$ws.Range("J73").Value = 4787.375
$ws.Range("K73").Value = 4745
$ws.Range("L73").Value = 4787.375
$ws.Range("M73").Value = -3809
$ws.Range("N73").Value = -6659.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 3646.6
$ws.Range("J107").Value = 3944.3333
$ws.Range("L107").Value = 3944.3333
$ws.Range("N107").Value = -7784.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1906
$ws.Range("I122").Value = 1990.25
$ws.Range("J122").Value = 1737.5
$ws.Range("K122").Value = 5970.75
$ws.Range("L122").Value = 5212.5
$ws.Range("M122").Value = -3520.75
$ws.Range("N122").Value = -10112.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1157002.9
$ws.Range("I122").Value = 1963784.8
$ws.Range("J122").Value = 4457.143
$ws.Range("K122").Value = 5891354.4
$ws.Range("L122").Value = 13371.429
$ws.Range("M122").Value = -5888904.4
$ws.Range("N122").Value = -18271.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3709.65
$ws.Range("I132").Value = 3154.5454
$ws.Range("J132").Value = 4388.1113
$ws.Range("K132").Value = 9463.636200000001
$ws.Range("L132").Value = 13164.3339
$ws.Range("M132").Value = -6933.636200000001
$ws.Range("N132").Value = -18224.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2233.3333
$ws.Range("J107").Value = 3000
$ws.Range("L107").Value = 9000
$ws.Range("N107").Value = -12840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2160.4443
$ws.Range("I122").Value = 1857.1428
$ws.Range("K122").Value = 5571.428400000001
$ws.Range("M122").Value = -3121.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 39464.5
$ws.Range("J123").Value = 39464.5
$ws.Range("L123").Value = 39464.5
$ws.Range("N123").Value = -49264.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 100429
$ws.Range("J138").Value = 100429
$ws.Range("L138").Value = 100429
$ws.Range("N138").Value = -110709
